$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22
# Data rows were re-shuffled between rows 22-28 and 45-56 (Id, TaxonId,
# Artnamn, Vetenskapligt namn, Auktor, Lokalnamn, Ost, Nord, Rapportör,
# Observatörer columns), per the upstream "Automatic update of files."
# commit. Each block below rewrites one row to its new content.
$ws.Range("A22").Value = 112323405
$ws.Range("B22").Value = 77650
$ws.Range("Q22").Value = 376903
$ws.Range("R22").Value = 6700268

# Row 23
$ws.Range("A23").Value = 112323404
$ws.Range("B23").Value = 78242
$ws.Range("E23").Value = 6453
$ws.Range("F23").Value = "Vedskivlav"
$ws.Range("G23").Value = "Hertelidea botryosa"
$ws.Range("H23").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q23").Value = 376956
$ws.Range("R23").Value = 6700282

# Row 24
$ws.Range("A24").Value = 112323406
$ws.Range("B24").Value = 77650
$ws.Range("Q24").Value = 376867
$ws.Range("R24").Value = 6700216

# Row 25
$ws.Range("B25").Value = 90837

# Row 26
$ws.Range("A26").Value = 112322582
$ws.Range("B26").Value = 77650
$ws.Range("E26").Value = 6425
$ws.Range("F26").Value = "Garnlav"
$ws.Range("G26").Value = "Alectoria sarmentosa"
$ws.Range("H26").Value = "(Ach.) Ach."
$ws.Range("P26").Value = "Gräsviggen, Vrm"
$ws.Range("Q26").Value = 376995
$ws.Range("R26").Value = 6700394
$ws.Range("AW26").Value = "Helena Malmestrand"
$ws.Range("AX26").Value = "Helena Malmestrand"

# Row 27
$ws.Range("A27").Value = 112323408
$ws.Range("B27").Value = 90837
$ws.Range("E27").Value = 5966
$ws.Range("F27").Value = "Motaggsvamp"
$ws.Range("G27").Value = "Sarcodon squamosus"
$ws.Range("H27").Value = "(Schaeff.) Quél."
$ws.Range("P27").Value = "Jonsmyren, Vrm"
$ws.Range("Q27").Value = 376863
$ws.Range("R27").Value = 6700246
$ws.Range("AW27").Value = "anders tedeholm"
$ws.Range("AX27").Value = "anders tedeholm"

# Row 28
$ws.Range("A28").Value = 112323403
$ws.Range("B28").Value = 77650
$ws.Range("E28").Value = 6425
$ws.Range("F28").Value = "Garnlav"
$ws.Range("G28").Value = "Alectoria sarmentosa"
$ws.Range("H28").Value = "(Ach.) Ach."
$ws.Range("Q28").Value = 376976
$ws.Range("R28").Value = 6700256

# Row 45
$ws.Range("A45").Value = 112323397
$ws.Range("B45").Value = 77650
$ws.Range("Q45").Value = 377010
$ws.Range("R45").Value = 6699884

# Row 46
$ws.Range("A46").Value = 112323396
$ws.Range("B46").Value = 77650
$ws.Range("P46").Value = "Jonsmyren, Vrm"
$ws.Range("Q46").Value = 376970
$ws.Range("R46").Value = 6699876
$ws.Range("AW46").Value = "anders tedeholm"
$ws.Range("AX46").Value = "anders tedeholm"

# Row 47
$ws.Range("A47").Value = 112322639
$ws.Range("B47").Value = 77650
$ws.Range("P47").Value = "Gräsviggen, Vrm"
$ws.Range("Q47").Value = 376970
$ws.Range("R47").Value = 6699950
$ws.Range("AW47").Value = "Helena Malmestrand"
$ws.Range("AX47").Value = "Helena Malmestrand"

# Row 48
$ws.Range("A48").Value = 112323395
$ws.Range("B48").Value = 90814
$ws.Range("D48").Value = "LC"
$ws.Range("E48").Value = 4364
$ws.Range("F48").Value = "Dropptaggsvamp"
$ws.Range("G48").Value = "Hydnellum ferrugineum"
$ws.Range("H48").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q48").Value = 376769
$ws.Range("R48").Value = 6699865

# Row 49
$ws.Range("A49").Value = 112322635
$ws.Range("B49").Value = 78242
$ws.Range("P49").Value = "Gräsviggen, Vrm"
$ws.Range("Q49").Value = 376947
$ws.Range("R49").Value = 6699898
$ws.Range("AW49").Value = "Helena Malmestrand"
$ws.Range("AX49").Value = "Helena Malmestrand"

# Row 50
$ws.Range("A50").Value = 112322604
$ws.Range("B50").Value = 77650
$ws.Range("D50").Value = "NT"
$ws.Range("E50").Value = 6425
$ws.Range("F50").Value = "Garnlav"
$ws.Range("G50").Value = "Alectoria sarmentosa"
$ws.Range("H50").Value = "(Ach.) Ach."
$ws.Range("P50").Value = "Gräsviggen, Vrm"
$ws.Range("Q50").Value = 377037
$ws.Range("R50").Value = 6699915
$ws.Range("AW50").Value = "Helena Malmestrand"
$ws.Range("AX50").Value = "Helena Malmestrand"

# Row 51
$ws.Range("A51").Value = 112323401
$ws.Range("B51").Value = 78242
$ws.Range("E51").Value = 6453
$ws.Range("F51").Value = "Vedskivlav"
$ws.Range("G51").Value = "Hertelidea botryosa"
$ws.Range("H51").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("P51").Value = "Jonsmyren, Vrm"
$ws.Range("Q51").Value = 376945
$ws.Range("R51").Value = 6700094
$ws.Range("AW51").Value = "anders tedeholm"
$ws.Range("AX51").Value = "anders tedeholm"

# Row 52
$ws.Range("A52").Value = 112323400
$ws.Range("B52").Value = 77650
$ws.Range("E52").Value = 6425
$ws.Range("F52").Value = "Garnlav"
$ws.Range("G52").Value = "Alectoria sarmentosa"
$ws.Range("H52").Value = "(Ach.) Ach."
$ws.Range("P52").Value = "Jonsmyren, Vrm"
$ws.Range("Q52").Value = 376916
$ws.Range("R52").Value = 6699968
$ws.Range("AW52").Value = "anders tedeholm"
$ws.Range("AX52").Value = "anders tedeholm"

# Row 53
$ws.Range("A53").Value = 112323399
$ws.Range("B53").Value = 77650
$ws.Range("E53").Value = 6425
$ws.Range("F53").Value = "Garnlav"
$ws.Range("G53").Value = "Alectoria sarmentosa"
$ws.Range("H53").Value = "(Ach.) Ach."
$ws.Range("P53").Value = "Jonsmyren, Vrm"
$ws.Range("Q53").Value = 377047
$ws.Range("R53").Value = 6699906
$ws.Range("AW53").Value = "anders tedeholm"
$ws.Range("AX53").Value = "anders tedeholm"

# Row 54
$ws.Range("A54").Value = 112323398
$ws.Range("B54").Value = 77650
$ws.Range("Q54").Value = 377056
$ws.Range("R54").Value = 6699887

# Row 55
$ws.Range("B55").Value = 78242

# Row 56
$ws.Range("A56").Value = 112322577
$ws.Range("B56").Value = 77403
$ws.Range("E56").Value = 228912
$ws.Range("F56").Value = "Mörk kolflarnlav"
$ws.Range("G56").Value = "Carbonicola myrmecina"
$ws.Range("H56").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("P56").Value = "Gräsviggen, Vrm"
$ws.Range("Q56").Value = 376932
$ws.Range("R56").Value = 6700074
$ws.Range("AW56").Value = "Helena Malmestrand"
$ws.Range("AX56").Value = "Helena Malmestrand"
Write-Output "Applied row re-shuffle across rows 22-28 and 45-56."
